# Ajout d'un système de grille externe ainsi qu'une sélection de grille aléatoire
# Append two new journal entries (rows 51 and 52) to the "Tableau1" table on
# the single worksheet, matching the existing formatting/formula pattern used
# by the rows above them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Grow the table (and its AutoFilter) from E5:M50 to E5:M52.
$lo.Resize($ws.Range("E5:M52"))

# Clone the formatting (number formats, alignment, wrap) of the last existing
# data row onto the two freshly added rows before filling in their values.
$ws.Range("E50:M50").Copy()
$ws.Range("E51:M51").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("E50:M50").Copy()
$ws.Range("E52:M52").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

$durationFormula = "=IF(ISBLANK(Tableau1[[#This Row],[Heure Début]]),`"`",Tableau1[[#This Row],[Heure fin]]-Tableau1[[#This Row],[Heure Début]])"

# --- Row 51 --------------------------------------------------------------
$ws.Cells.Item(51, 5).Value = 44287
$ws.Cells.Item(51, 6).Value = 0.44444444444444442
$ws.Cells.Item(51, 7).Value = 0.50347222222222221
$ws.Cells.Item(51, 8).Formula = $durationFormula
$ws.Cells.Item(51, 9).Value = "Développement"
$ws.Cells.Item(51, 10).Value = "Implémenter des fichier stoquant les grilles"
$ws.Cells.Item(51, 11).Value = "CPNV"
$ws.Cells.Item(51, 12).Value = "Faire que le programme choisissent au hazard une grille pour la partie"
$ws.Cells.Item(51, 13).Value = "VWM`nhttps://www.programmingsimplified.com/c-program-generate-random-numbers"
$ws.Rows.Item(51).RowHeight = 57.6

# --- Row 52 --------------------------------------------------------------
$ws.Cells.Item(52, 5).Value = 44287
$ws.Cells.Item(52, 6).Value = 0.50347222222222221
$ws.Cells.Item(52, 7).Value = 0.51041666666666663
$ws.Cells.Item(52, 8).Formula = $durationFormula
$ws.Cells.Item(52, 9).Value = "Développement"
$ws.Cells.Item(52, 10).Value = "Correction de warning dans le code"
$ws.Cells.Item(52, 11).Value = "CPNV"
$ws.Cells.Item(52, 12).Value = "Réécrire le code pour qu'il n'y ait plus d'érreur"
$ws.Rows.Item(52).RowHeight = 43.2

# Match the author's final cursor position after typing the new rows.
$ws.Range("H53").Select()
